# Resort the worksheet tabs: "总计" should come before "2021-Q1".
# (The data/content of each named sheet is unchanged - only the tab
# order changes. "2021-Q1" stays the active/selected sheet.)
$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item("总计")
$firstSheet = $wb.Worksheets.Item(1)

# Move "总计" in front of whatever sheet is currently first ("2021-Q1").
$totalSheet.Move($firstSheet)

# Keep "2021-Q1" as the active/selected sheet, as it was before the reorder.
$q1Sheet = $wb.Worksheets.Item("2021-Q1")
$q1Sheet.Activate()
